# Update Daily Report: 2026-02-10
# Adds the next business day's depository data (date serial 46062 = 2026-02-09)
# to Daily_Data, and refreshes the two summary sheets (Today_Summary,
# Monthly_Stats) so their latest-day / month-to-date figures reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append rows 552-573 for date 46062
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(46062, "ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @(46062, "ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(46062, "BRINK'S, INC. Registered", 75623.302, 0, 0, 0, -2268.519, 73354.783),
    @(46062, "BRINK'S, INC. Eligible", 83553.32799999999, 0, 0, 0, 2268.519, 85821.84699999999),
    @(46062, "CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @(46062, "CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @(46062, "DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @(46062, "DELAWARE DEPOSITORY Eligible", 18459.584, 0, 0, 0, 0, 18459.584),
    @(46062, "HSBC BANK, USA Registered", 1394.758, 0, 0, 0, 0, 1394.758),
    @(46062, "HSBC BANK, USA Eligible", 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @(46062, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @(46062, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @(46062, "JP MORGAN CHASE BANK NA Registered", 114985.579, 0, 0, 0, -924.158, 114061.421),
    @(46062, "JP MORGAN CHASE BANK NA Eligible", 75484.511, 0, 0, 0, 924.158, 76408.66899999999),
    @(46062, "LOOMIS INTERNATIONAL (US) LLC Registered", 63745.991, 0, 0, 0, -2588.547, 61157.444),
    @(46062, "LOOMIS INTERNATIONAL (US) LLC Eligible", 69005.64, 0, 0, 0, 2588.547, 71594.18700000001),
    @(46062, "MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @(46062, "MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(46062, "MANFRA, TORDELLA & BROOKES, LLC Registered", 50220.42, 0, 0, 0, -300.172, 49920.248),
    @(46062, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 1804.683, 0, 0, 0, 300.172, 2104.855),
    @(46062, "STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @(46062, "STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

$startRow = 552
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $row[0]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Today_Summary: refresh Eligible/Registered split for the
#    depositories that moved metal between categories today
#    (Total_Stock per depository is unchanged).
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("Today_Summary")

$today.Range("B3").Value = 85821.84699999999   # BRINK'S, INC. Eligible
$today.Range("C3").Value = 73354.783           # BRINK'S, INC. Registered

$today.Range("B8").Value = 76408.66899999999   # JP MORGAN CHASE BANK NA Eligible
$today.Range("C8").Value = 114061.421          # JP MORGAN CHASE BANK NA Registered

$today.Range("B9").Value = 71594.18700000001   # LOOMIS INTERNATIONAL (US) LLC Eligible
$today.Range("C9").Value = 61157.444           # LOOMIS INTERNATIONAL (US) LLC Registered

$today.Range("B11").Value = 2104.855           # MANFRA, TORDELLA & BROOKES, LLC Eligible
$today.Range("C11").Value = 49920.248          # MANFRA, TORDELLA & BROOKES, LLC Registered

# ---------------------------------------------------------------------
# 3) Monthly_Stats: refresh the 2026-02 Eligible/Registered column
#    totals and the detail rows for the same four depositories.
# ---------------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly_Stats")

$monthly.Range("B2").Value = 263687.196        # 2026-02 Eligible total
$monthly.Range("C2").Value = 319682.013        # 2026-02 Registered total

$monthly.Range("E10").Value = 85821.84699999999  # BRINK'S, INC. Eligible
$monthly.Range("E11").Value = 73354.783          # BRINK'S, INC. Registered

$monthly.Range("E20").Value = 76408.66899999999  # JP MORGAN CHASE BANK NA Eligible
$monthly.Range("E21").Value = 114061.421         # JP MORGAN CHASE BANK NA Registered

$monthly.Range("E22").Value = 71594.18700000001  # LOOMIS INTERNATIONAL (US) LLC Eligible
$monthly.Range("E23").Value = 61157.444          # LOOMIS INTERNATIONAL (US) LLC Registered

$monthly.Range("E26").Value = 2104.855           # MANFRA, TORDELLA & BROOKES, LLC Eligible
$monthly.Range("E27").Value = 49920.248          # MANFRA, TORDELLA & BROOKES, LLC Registered
